# TC01_CDS_Filter_Study-GECCO-OICR.xlsx
# "Filter - Study - Test Suit"
#
# The startup sheet lists the tabs that get generated/filtered for this
# test case. Rename the "CasesTab" row to "ParticipantsTab" (its query /
# file-name columns already describe the participant-level query), and
# leave the cursor on that renamed cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$ws.Range("A2").Value = "ParticipantsTab"
$ws.Range("A2").Select()
